# Update the "ACHIEVEMENTS" (Class Time Tracker) worksheet:
#  - Enter actual-time values in column D for rows 14-18
#  - Enter estimated-time values in column C for rows 20-23
#  - Update the view's scroll position / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACHIEVEMENTS")

# Column D (Actual Time) entries for rows 14-18
$ws.Range("D14").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("D18").Value = 2

# Column C (Estimated Time) entries for rows 20-23
$ws.Range("C20").Value = 3
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 2
$ws.Range("C23").Value = 2

# Move the view so row 14 is at the top and select D19
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("D19").Select()
